$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5822
$ws.Range("I3").Value = 7493
$ws.Range("J3").Value = 6219
$ws.Range("C4").Value = 1836
$ws.Range("D4").Value = 1961
$ws.Range("E4").Value = 2009
$ws.Range("J4").Value = 1345
$ws.Range("J5").Value = 475
$ws.Range("J6").Value = 7939
$ws.Range("C7").Value = 28380
$ws.Range("D7").Value = 28151
$ws.Range("E7").Value = 26015
$ws.Range("I7").Value = 26229
$ws.Range("J7").Value = 21800

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 197
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 419
$ws.Range("J5").Value = 36
$ws.Range("J6").Value = 463
$ws.Range("J7").Value = 1370

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 134
$ws.Range("J3").Value = 162
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 443

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 243
$ws.Range("J3").Value = 338
$ws.Range("J6").Value = 347
$ws.Range("J7").Value = 1012

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 113
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 227
$ws.Range("J7").Value = 670

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 97
$ws.Range("J7").Value = 342

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 96
$ws.Range("J5").Value = 69
$ws.Range("J7").Value = 644
$ws.Range("J8").Value = 1370
$ws.Range("J9").Value = 108
$ws.Range("J11").Value = 341
$ws.Range("J15").Value = 241
$ws.Range("J18").Value = 181
$ws.Range("J19").Value = 642
$ws.Range("J20").Value = 450
$ws.Range("J23").Value = 206
$ws.Range("J25").Value = 109
$ws.Range("J29").Value = 1211
$ws.Range("J31").Value = 197
$ws.Range("J33").Value = 1012
$ws.Range("J37").Value = 670
$ws.Range("J41").Value = 140
$ws.Range("J42").Value = 915
$ws.Range("J44").Value = 165
$ws.Range("J47").Value = 166
$ws.Range("J48").Value = 259
$ws.Range("J52").Value = 547
$ws.Range("J53").Value = 301
$ws.Range("J54").Value = 425
$ws.Range("J55").Value = 297
$ws.Range("J60").Value = 130
$ws.Range("C63").Value = 266
$ws.Range("D63").Value = 345
$ws.Range("E63").Value = 352
$ws.Range("I63").Value = 242
$ws.Range("J67").Value = 823
$ws.Range("J73").Value = 209
$ws.Range("J75").Value = 65
$ws.Range("J76").Value = 330
$ws.Range("J78").Value = 271
$ws.Range("J79").Value = 622
$ws.Range("J83").Value = 443
$ws.Range("J85").Value = 902
$ws.Range("J86").Value = 137
$ws.Range("J90").Value = 233
$ws.Range("J91").Value = 246
$ws.Range("J93").Value = 96
$ws.Range("J94").Value = 222
$ws.Range("J95").Value = 323
$ws.Range("J98").Value = 155
$ws.Range("J99").Value = 342
$ws.Range("J100").Value = 40
$ws.Range("C101").Value = 28380
$ws.Range("D101").Value = 28151
$ws.Range("E101").Value = 26015
$ws.Range("I101").Value = 26229
$ws.Range("J101").Value = 21800

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 312
$ws.Range("J6").Value = 220
$ws.Range("J7").Value = 823

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 101
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 425

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 425
$ws.Range("J6").Value = 311
$ws.Range("J7").Value = 1211

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J6").Value = 245
$ws.Range("J7").Value = 642

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 50
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 53
$ws.Range("J7").Value = 330

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 198
$ws.Range("J3").Value = 187
$ws.Range("J6").Value = 472
$ws.Range("J7").Value = 915

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 73
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J5").Value = 7
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 101
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 246

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 175
$ws.Range("J7").Value = 622

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 155
$ws.Range("J7").Value = 450

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 197
$ws.Range("J3").Value = 194
$ws.Range("J4").Value = 27
$ws.Range("J6").Value = 207
$ws.Range("J7").Value = 644

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 43
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 241

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 68
$ws.Range("J6").Value = 142
$ws.Range("J7").Value = 341

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 73
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 238
$ws.Range("J4").Value = 59
$ws.Range("J7").Value = 902

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J4").Value = 20
$ws.Range("J6").Value = 221
$ws.Range("J7").Value = 547

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 32
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 96
